# Applies the "Updated cryptos list" GitHub Actions refresh to Sheet1.
# Each data cell is written with a leading apostrophe so Excel keeps the
# value as text (prices/volumes are formatted strings, not numbers), then
# ClearFormats() drops the transient "quote prefix" flag so the cell style
# is left exactly as it was before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# Row 2
Set-TextValue "D2" "27.039.16"
Set-TextValue "E2" "  +0.58%  "
# Row 3
Set-TextValue "D3" "1.683.04"
Set-TextValue "E3" "  +0.82%  "
# Row 4
Set-TextValue "E4" "  +0.02%  "
# Row 5
Set-TextValue "D5" "216.14"
Set-TextValue "E5" "  +0.23%  "
# Row 6
Set-TextValue "E6" "  -2.35%  "
# Row 7
Set-TextValue "E7" "  -0.04%  "
# Row 8
Set-TextValue "D8" "21.60"
Set-TextValue "E8" "  +6.57%  "
# Row 9
Set-TextValue "E9" "  -0.23%  "
# Row 10
Set-TextValue "D10" "0.0623"
Set-TextValue "E10" "  +0.66%  "
# Row 11
Set-TextValue "D11" "0.0891"
Set-TextValue "E11" "  -0.24%  "
# Row 12
Set-TextValue "D12" "1.919.55"
Set-TextValue "E12" "  +0.79%  "
# Row 13
Set-TextValue "D13" "1.678.14"
Set-TextValue "E13" "  +0.40%  "
# Row 14
Set-TextValue "E14" "  +0.49%  "
# Row 15
Set-TextValue "E15" "  +1.81%  "
# Row 16
Set-TextValue "D16" "66.33"
Set-TextValue "E16" "  +0.83%  "
# Row 17
Set-TextValue "B17" "WrappedBTC"
Set-TextValue "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "27.043.55"
Set-TextValue "E17" "  +0.52%  "
# Row 18
Set-TextValue "B18" "Chainlink"
Set-TextValue "C18" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D18" "8.19"
Set-TextValue "E18" "  +5.19%  "
# Row 19
Set-TextValue "D19" "236.68"
Set-TextValue "E19" "  +2.04%  "
# Row 20
Set-TextValue "E20" "  +0.56%  "
# Row 21
Set-TextValue "E21" "  +0.03%  "
# Row 22
Set-TextValue "E22" "  +0.00%  "
# Row 23
Set-TextValue "D23" "9.28"
Set-TextValue "E23" "  +0.88%  "
# Row 24
Set-TextValue "E24" "  -4.11%  "
# Row 25
Set-TextValue "D25" "147.06"
Set-TextValue "E25" "  +1.02%  "
# Row 26
Set-TextValue "D26" "16.80"
Set-TextValue "E26" "  +5.58%  "
# Row 27
Set-TextValue "D27" "7.25"
Set-TextValue "E27" "  +1.70%  "
# Row 28
Set-TextValue "E28" "  -2.94%  "
# Row 29
Set-TextValue "E29" "  -0.11%  "
# Row 30
Set-TextValue "E30" "  +0.47%  "
# Row 31
Set-TextValue "D31" "1.17"
Set-TextValue "E31" "  -0.21%  "
# Row 32
Set-TextValue "D32" "3.35"
Set-TextValue "E32" "  +0.30%  "
# Row 33
Set-TextValue "D33" "1.522.91"
Set-TextValue "E33" "  +3.96%  "
# Row 34
Set-TextValue "E34" "  +0.76%  "
# Row 35
Set-TextValue "E35" "  +4.62%  "
# Row 36
Set-TextValue "E36" "  -0.40%  "
# Row 37
Set-TextValue "D37" "0.590"
Set-TextValue "E37" "  +3.30%  "
# Row 38
Set-TextValue "B38" "ARBITRUM"
Set-TextValue "C38" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D38" "0.921"
Set-TextValue "E38" "  +2.43%  "
# Row 39
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0175"
Set-TextValue "E39" "  +3.86%  "
# Row 40
Set-TextValue "E40" "  +6.70%  "
# Row 41
Set-TextValue "E41" "  -0.91%  "
# Row 42
Set-TextValue "E42" "  +0.01%  "
# Row 43
Set-TextValue "D43" "68.12"
Set-TextValue "E43" "  +3.62%  "
# Row 44
Set-TextValue "D44" "2.27"
Set-TextValue "E44" "  -0.54%  "
# Row 45
Set-TextValue "D45" "1.823.49"
Set-TextValue "E45" "  +0.40%  "
# Row 46
Set-TextValue "E46" "  +0.19%  "
# Row 47
Set-TextValue "D47" "90.24"
# Row 48
Set-TextValue "B48" "BabyDogeCoin"
Set-TextValue "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.0₆0105"
Set-TextValue "E48" "  +0.12%  "
# Row 49
Set-TextValue "B49" "Algorand"
Set-TextValue "C49" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D49" "0.104"
Set-TextValue "E49" "  +4.22%  "
# Row 50
Set-TextValue "B50" "RenderToken"
Set-TextValue "C50" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "1.52"
Set-TextValue "E50" "  -0.39%  "
# Row 51
Set-TextValue "B51" "EnergySwap"
Set-TextValue "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.94"
Set-TextValue "E51" "  +4.64%  "
